$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 4 notes: clarify that dataset excluded non-U.S. companies
# AND employees, and that the focus is U.S. companies/residents.
# ---------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$notes4 = $s4.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notes4.Paragraphs(1,1).Text = "This dataset contained information about companies and employees that were not located in the U.S. For the purpose of my analysis, I excluded those and only focused on U.S. based companies and U.S. residents. "

# ---------------------------------------------------------------
# Slide 9 notes: "medium companies" -> "medium sized companies"
# ---------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$notes9 = $s9.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notes9.Paragraphs(9,1).Text = "Surprisingly, medium sized companies had the highest amount of in person jobs and very little hybrid work offered. "

# ---------------------------------------------------------------
# Slide 11 notes: rewrite the three "future project plans" notes
# ---------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$notes11 = $s11.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notes11.Paragraphs(1,1).Text = "If I had more time, I" + [char]8217 + "d like to explore the different pay based on employment types like part-time, contracting, and freelance."
$notes11.Paragraphs(2,1).Text = "I would also like to compare remote work salaries to in person and hybrid salaries."
$notes11.Paragraphs(3,1).Text = "It would also be interesting to compare U.S. average pay to other countries."

# ---------------------------------------------------------------
# Slide 11 body textbox: rewrite the three matching bullet lines
# ---------------------------------------------------------------
$bullets = $s11.Shapes.Item(3).TextFrame.TextRange
$bullets.Paragraphs(1,1).Text = "Pay based on type of employment"
$bullets.Paragraphs(2,1).Text = "Remote work compared to salary"
$bullets.Paragraphs(3,1).Text = "U.S. average pay compared to other countries"
